# Add two new Mac-Addresses (10 new detail rows for two new machines)
# to the reg_center_machine_device master data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append: regcntr_id, machine_id, device_id
# then fixed columns D=eng, E=TRUE, F=superadmin, G=now()
$newRows = @(
    @(10001, 10030, 3000166),
    @(10001, 10030, 3000167),
    @(10001, 10030, 3000168),
    @(10001, 10030, 3000169),
    @(10001, 10030, 3000170),
    @(10001, 10031, 3000171),
    @(10001, 10031, 3000172),
    @(10001, 10031, 3000173),
    @(10001, 10031, 3000174),
    @(10001, 10031, 3000175)
)

$startRow = 147

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = "now()"
}

$lastRow = $startRow + $newRows.Count - 1

# Update the selection to reflect where the editor ended up after adding rows
$ws.Range("H149").Select()
